# Updated as of March 23 18:00.
# Appends the new observation (data_collection="live", source="govt_canada_ph",
# date 2020-03-23 18:00) to both the "longform" (per-province wide table) and
# "shortform" (province/case_type long table) sheets, then restores each
# sheet's selection to the newly-appended cells.

$wb = $excel.ActiveWorkbook
$wsLong  = $wb.Worksheets.Item("longform")
$wsShort = $wb.Worksheets.Item("shortform")

# ---------------------------------------------------------------------------
# 1. longform (sheet2): new row 31
# ---------------------------------------------------------------------------
$longRow = 31

$wsLong.Cells.Item($longRow, 1).Value = "live"
$wsLong.Cells.Item($longRow, 2).Value = "govt_canada_ph"

# Date cell: copy the format (style) of the cell directly above so the new
# cell reuses the existing date-time style instead of minting a new one.
$wsLong.Cells.Item($longRow, 3).Value = 43913.75
$wsLong.Cells.Item($longRow - 1, 3).Copy()
$wsLong.Cells.Item($longRow, 3).PasteSpecial(-4122)

$longVals = @(472,0,13,301,0,1,65,1,0,11,9,0,503,0,6,221,407,4,4,20,0,9,8,0,41,0,0,3,0,0,13,0,0,2,0,0,1,0,0)
for ($i = 0; $i -lt $longVals.Length; $i++) {
    $wsLong.Cells.Item($longRow, 4 + $i).Value = $longVals[$i]
}

# ---------------------------------------------------------------------------
# 2. shortform (sheet3): new rows 288-326 (13 provinces x 3 case types)
# ---------------------------------------------------------------------------
$provinces = @("BC","AB","SK","MB","ON","QC","NL","NB","NS","PEI","Repat","YK","NT")
$caseTypes = @("conf","prob","deaths")
$shortVals = @(
    472,0,13,
    301,0,1,
    65,1,0,
    11,9,0,
    503,0,6,
    221,407,4,
    4,20,0,
    9,8,0,
    41,0,0,
    3,0,0,
    13,0,0,
    2,0,0,
    1,0,0
)

$shortRow = 288
$k = 0
for ($p = 0; $p -lt $provinces.Length; $p++) {
    for ($c = 0; $c -lt $caseTypes.Length; $c++) {
        $wsShort.Cells.Item($shortRow, 1).Value = "live"
        $wsShort.Cells.Item($shortRow, 2).Value = "govt_canada_ph"

        $wsShort.Cells.Item($shortRow, 3).Value = 43913.75
        $wsShort.Cells.Item($shortRow - 1, 3).Copy()
        $wsShort.Cells.Item($shortRow, 3).PasteSpecial(-4122)

        $wsShort.Cells.Item($shortRow, 4).Value = $provinces[$p]
        $wsShort.Cells.Item($shortRow, 5).Value = $caseTypes[$c]
        $wsShort.Cells.Item($shortRow, 6).Value = $shortVals[$k]

        $shortRow = $shortRow + 1
        $k = $k + 1
    }
}

# ---------------------------------------------------------------------------
# 3. Restore selections to the newly-appended cells on each sheet.
#    (Range.Select requires the sheet to be active first; re-activate
#    "shortform" afterwards since that is the workbook's original active tab.)
# ---------------------------------------------------------------------------
$wsLong.Activate()
$wsLong.Range("C31").Select()

$wsShort.Activate()
$wsShort.Range("C288:C326").Select()
